$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.106.45"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "3.130.14"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.32%  "

$ws.Range("D6").Value = "137.13"
$ws.Range("E6").Value = "  -4.50%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "3.123.25"
$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("E9").Value = "  -1.69%  "

$ws.Range("E10").Value = "  -2.85%  "

$ws.Range("E11").Value = "  -1.59%  "

$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  -3.11%  "

$ws.Range("D14").Value = "34.42"
$ws.Range("E14").Value = "  -3.22%  "

$ws.Range("D15").Value = "3.639.85"
$ws.Range("E15").Value = "  -0.38%  "

$ws.Range("E16").Value = "  +1.69%  "

$ws.Range("D17").Value = "63.128.75"
$ws.Range("E17").Value = "  -1.83%  "

$ws.Range("D18").Value = "3.127.20"
$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("D19").Value = "6.76"
$ws.Range("E19").Value = "  -1.41%  "

$ws.Range("D20").Value = "477.08"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").Value = "14.21"
$ws.Range("E21").Value = "  -3.41%  "

$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  -2.95%  "

$ws.Range("D23").Value = "7.71"
$ws.Range("E23").Value = "  -1.26%  "

$ws.Range("D24").Value = "87.35"
$ws.Range("E24").Value = "  +2.55%  "

$ws.Range("D25").Value = "13.08"
$ws.Range("E25").Value = "  -4.10%  "

$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("E27").Value = "  -2.04%  "

$ws.Range("D28").Value = "7.22"
$ws.Range("E28").Value = "  -2.13%  "

$ws.Range("E29").Value = "  -6.51%  "

$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.68%  "

$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("E33").Value = "  -7.54%  "

$ws.Range("E34").Value = "  -3.31%  "

$ws.Range("E35").Value = "  -2.55%  "

$ws.Range("E36").Value = "  -1.66%  "

$ws.Range("D37").Value = "51.93"

$ws.Range("E38").Value = "  -4.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0390"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.43%  "

$ws.Range("D40").Value = "422.02"
$ws.Range("E40").Value = "  -7.25%  "

$ws.Range("E41").Value = "  -0.75%  "

$ws.Range("E42").Value = "  -0.45%  "

$ws.Range("E43").Value = "  -10.57%  "

$ws.Range("D44").Value = "2.891.59"
$ws.Range("E44").Value = "  +1.02%  "

$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").Value = "2.13"
$ws.Range("E46").Value = "  -5.82%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").Value = "25.85"
$ws.Range("E48").Value = "  -2.24%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("E50").Value = "  -5.77%  "

$ws.Range("D51").Value = "118.78"
$ws.Range("E51").Value = "  -1.35%  "
